$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.077.65'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '3.485.76'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '594.79'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '182.54'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  +3.39%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '3.480.93'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '4.082.82'
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.32'
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('E15').Value = '  -0.87%  '
$ws.Range('D16').Value = '68.060.18'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').Value = '3.487.18'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.23'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.14'
$ws.Range('E20').Value = '  -4.01%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '394.92'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.96'
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('E23').Value = '  +2.41%  '
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '72.24'
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.46'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.41'
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.65'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.38'
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  -5.89%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '161.74'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.891'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.84'
$ws.Range('E40').Value = '  +5.43%  '
$ws.Range('E41').Value = '  -3.77%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.68'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.77'
$ws.Range('E43').Value = '  -4.63%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '26.30'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('D46').Value = '2.750.62'
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('E47').Value = '  -5.47%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '41.57'
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '332.08'
$ws.Range('E50').Value = '  -3.86%  '
$ws.Range('E51').Value = '  -3.75%  '
